$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the columns we will touch so Excel keeps values as text
# (matches original inlineStr / text-typed cells), then restore default style after.
$ws.Range("B47:C48").NumberFormat = "@"
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "261.27"
$ws.Range("E2").Value = "1.85%"
$ws.Range("G2").Value = "10"
$ws.Range("D3").Value = "27.43"
$ws.Range("E3").Value = "1.85%"
$ws.Range("G3").Value = "10"
$ws.Range("D4").Value = "4.697"
$ws.Range("E4").Value = "-0.17%"
$ws.Range("G4").Value = "10"
$ws.Range("D5").Value = "0.06091"
$ws.Range("E5").Value = "2.79%"
$ws.Range("G5").Value = "10"
$ws.Range("D6").Value = "6.679"
$ws.Range("E6").Value = "1.11%"
$ws.Range("G6").Value = "10"
$ws.Range("D7").Value = "0.8461"
$ws.Range("E7").Value = "-0.51%"
$ws.Range("G7").Value = "10"
$ws.Range("D8").Value = "0.9270"
$ws.Range("E8").Value = "0.53%"
$ws.Range("G8").Value = "10"
$ws.Range("D9").Value = "0.1405"
$ws.Range("E9").Value = "1.61%"
$ws.Range("G9").Value = "10"
$ws.Range("D10").Value = "0.04695"
$ws.Range("E10").Value = "12.87%"
$ws.Range("G10").Value = "10"
$ws.Range("D11").Value = "0.07109"
$ws.Range("E11").Value = "1.46%"
$ws.Range("G11").Value = "10"
$ws.Range("D12").Value = "0.03092"
$ws.Range("E12").Value = "1.31%"
$ws.Range("G12").Value = "10"
$ws.Range("D13").Value = "0.09069"
$ws.Range("E13").Value = "-0.35%"
$ws.Range("G13").Value = "10"
$ws.Range("D14").Value = "0.001533"
$ws.Range("E14").Value = "-0.59%"
$ws.Range("G14").Value = "10"
$ws.Range("D15").Value = "0.0006069"
$ws.Range("E15").Value = "-94.10%"
$ws.Range("G15").Value = "10"
$ws.Range("D16").Value = "0.006085"
$ws.Range("E16").Value = "-1.59%"
$ws.Range("G16").Value = "10"
$ws.Range("D17").Value = "3.449"
$ws.Range("E17").Value = "-0.60%"
$ws.Range("G17").Value = "10"
$ws.Range("D18").Value = "3.141"
$ws.Range("E18").Value = "-0.83%"
$ws.Range("G18").Value = "10"
$ws.Range("E19").Value = "-0.61%"
$ws.Range("G19").Value = "10"
$ws.Range("E20").Value = "2.27%"
$ws.Range("G20").Value = "10"
$ws.Range("D21").Value = "0.1305"
$ws.Range("E21").Value = "0.78%"
$ws.Range("G21").Value = "10"
$ws.Range("D22").Value = "4.089"
$ws.Range("E22").Value = "4.80%"
$ws.Range("G22").Value = "10"
$ws.Range("D23").Value = "0.04232"
$ws.Range("E23").Value = "-0.02%"
$ws.Range("G23").Value = "10"
$ws.Range("E24").Value = "0.32%"
$ws.Range("G24").Value = "10"
$ws.Range("D25").Value = "0.003786"
$ws.Range("E25").Value = "-11.44%"
$ws.Range("G25").Value = "10"
$ws.Range("G26").Value = "10"
$ws.Range("E27").Value = "3.46%"
$ws.Range("G27").Value = "10"
$ws.Range("G28").Value = "10"
$ws.Range("G29").Value = "10"
$ws.Range("G30").Value = "10"
$ws.Range("G31").Value = "10"
$ws.Range("G32").Value = "10"
$ws.Range("G33").Value = "10"
$ws.Range("G34").Value = "10"
$ws.Range("G35").Value = "10"
$ws.Range("G36").Value = "10"
$ws.Range("G37").Value = "10"
$ws.Range("G38").Value = "10"
$ws.Range("G39").Value = "10"
$ws.Range("E40").Value = "2.31%"
$ws.Range("G40").Value = "10"
$ws.Range("D41").Value = "0.1114"
$ws.Range("E41").Value = "1.37%"
$ws.Range("G41").Value = "10"
$ws.Range("D42").Value = "0.004091"
$ws.Range("E42").Value = "-34.29%"
$ws.Range("G42").Value = "10"
$ws.Range("D43").Value = "0.01630"
$ws.Range("E43").Value = "15.53%"
$ws.Range("G43").Value = "10"
$ws.Range("E44").Value = "16.74%"
$ws.Range("G44").Value = "10"
$ws.Range("D45").Value = "0.00005142"
$ws.Range("E45").Value = "-4.02%"
$ws.Range("G45").Value = "10"
$ws.Range("E46").Value = "0.09%"
$ws.Range("G46").Value = "10"
$ws.Range("B47").Value = "BOLO"
$ws.Range("C47").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D47").Value = "0.1358"
$ws.Range("E47").Value = "-45.67%"
$ws.Range("G47").Value = "10"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "0.05445"
$ws.Range("E48").Value = "18.18%"
$ws.Range("G48").Value = "10"
$ws.Range("E49").Value = "0.09%"
$ws.Range("G49").Value = "10"
$ws.Range("E50").Value = "0.09%"
$ws.Range("G50").Value = "10"
$ws.Range("G51").Value = "10"

# Restore default (Normal) style so formatting matches the original workbook
$ws.Range("B47:C48").Style = "Normal"
$ws.Range("D2:D51").Style = "Normal"
$ws.Range("E2:E51").Style = "Normal"
$ws.Range("G2:G51").Style = "Normal"
